$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pretty-printed JSON content that replaces the old Python-dict-repr string.
$newText = @'
questions = [
    {
        "title": "Which of the following attitudes is BEST suited for a career in the hospitality industry?",
        "ques_type": 2,
        "options": [
            "Cost consciousness",
            "Customer orientation",
            "Entrepreneurial spirit",
            "Environmental awareness"
        ],
        "score": "Customer orientation"
    },
    {
        "title": "Assume you are the floor manager walking around your restaurant during the busy lunch hour. You find the new waiter having trouble carrying hot plates of food from the kitchen to the service table.Which of the following is the BEST action that you can suggest?",
        "ques_type": 2,
        "options": [
            "Ask the kitchen staff to re-plate the dishes on cooler platters.",
            "Ask the waiter to bring in the food only after it is no longer hot.",
            "Ask the waiter to use service cloths or a tray to carry the item",
            "Give the waiter a different task."
        ],
        "score": "Ask the waiter to use service cloths or a tray to carry the item"
    },
    {
        "title": "Assume you are working at the front desk of a hotel. A guest who has a booking for a future date calls up to request certain additional paid facilities.What is the best way to end the call?",
        "ques_type": 2,
        "options": [
            "Tell the caller that you have everything you need, say goodbye or any other pleasantries, and hang up.",
            "Give a basic call summary, ask if there is anything else you can help out with, and finish with a courtesy statement.",
            "Let the caller know that you don\u2019t need any further information, thank them, and stay on the line until the customer decides to hang up.",
            "Ask the caller to repeat the requests one final time and re-confirm your understanding. Once that is verified, you may politely end the call. "
        ],
        "score": "Give a basic call summary, ask if there is anything else you can help out with, and finish with a courtesy statement."
    },
    {
        "title": "Which of the following is an example of cross-contamination of food?",
        "ques_type": 2,
        "options": [
            "Cooking poultry, lamb, and veal on the same restaurant grill.",
            "Serving vegetarian and non-vegetarian food in the same restaurant.",
            "Handling raw meat and vegetables without washing hands in between.",
            "Storing fish, dairy products, and meat in different refrigerator compartments."
        ],
        "score": "Handling raw meat and vegetables without washing hands in between."
    }
]
'@

# Write the long multi-line text into a scratch cell first and copy/paste it
# into A1 afterwards. Writing directly into A1 with embedded newlines makes
# the sheet auto-fit that row's height; routing the text through a
# copy/paste instead leaves row 1 at its normal (non-custom) height, which
# is what the target sheet looks like.
$ws.Range("Z1").Value = $newText
$ws.Range("Z1").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Drop the old row 1 (numeric "0" cell, bold/bordered style) and the now
# vacated original row 2, so the new text ends up alone on row 1 with the
# default (unstyled) format.
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()
